# Word COM-interop script implementing the "Better 3.2 and 3.3" edit.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Text edits (paragraph count unaffected at this point)
# ---------------------------------------------------------------------

# --- Question 3.2 answer paragraph (Word paragraph #45) ---
$p32 = $d.Paragraphs.Item(45)

$r = $p32.Range
[void]$r.Find.Execute(
    "values which they can represent", $false, $false, $false, $false,
    $false, $true, 1, $false, "values they can represent", 2)

$r = $p32.Range
[void]$r.Find.Execute(
    "both char and unsigned char", $false, $false, $false, $false,
    $false, $true, 1, $false, "both (signed) char and unsigned char", 2)

$r = $p32.Range
[void]$r.Find.Execute(
    "not affect the amount of information that can be stored within them",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "not affect anything", 2)

# --- Question 3.3 first answer paragraph (Word paragraph #53) ---
$p33a = $d.Paragraphs.Item(53)
$r = $p33a.Range
[void]$r.Find.Execute(
    "freed as it calls free_map which has the length parameter.",
    $false, $false, $false, $false, $false, $true, 1, $false, "freed.", 2)

# ---------------------------------------------------------------------
# 2) Paragraph alignment changes (still stable indices)
# ---------------------------------------------------------------------
$d.Paragraphs.Item(46).Range.ParagraphFormat.Alignment = 3   # blank after 3.2 answer
$d.Paragraphs.Item(51).Range.ParagraphFormat.Alignment = 3   # "Question 3.3 (1 mark)"
$d.Paragraphs.Item(52).Range.ParagraphFormat.Alignment = 3   # blank before 3.3 answer
$d.Paragraphs.Item(53).Range.ParagraphFormat.Alignment = 3   # "Yes, it needs to know..."
$d.Paragraphs.Item(54).Range.ParagraphFormat.Alignment = 3   # blank
$d.Paragraphs.Item(55).Range.ParagraphFormat.Alignment = 3   # "To get this information..."

# ---------------------------------------------------------------------
# 3) Insert the new "Besides, ..." paragraph right after paragraph 46
#    (the blank, bordered paragraph that follows the 3.2 answer).
# ---------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(46)
$anchor.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(47)
$newRange = $newPara.Range
$newRange.ParagraphFormat.Alignment = 3
$newRange.Text = "Besides, choosing to use char (signed or unsigned) in the 1st place is probably just because it is the smallest primitive data type in C and hence using it results in space savings. In reality, using other data types that can represent the 2 states of allocated and freed will work too."

# Make the "st" in "1st" superscript.
$newPara = $d.Paragraphs.Item(47)
$pStart = $newPara.Range.Start
$findRange = $d.Range($pStart, $newPara.Range.End)
[void]$findRange.Find.Execute("st place", $false, $false, $false, $false, $false, $true)
$stStart = $findRange.Start
$stRange = $d.Range($stStart, $stStart + 2)
$stRange.Font.Superscript = $true

# ---------------------------------------------------------------------
# 4) Remove 3 of the 4 blank "u=single" paragraphs that used to sit
#    between the 3.2 answer block and "Question 3.3". After the insert
#    above they now live at indices 48-51; keep one, delete the rest.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(48).Range.Delete()
$d.Paragraphs.Item(48).Range.Delete()
$d.Paragraphs.Item(48).Range.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
